$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Try simple value+style copy, bottom-up shift within column C only, rows 6..20
for ($r = 20; $r -ge 7; $r--) {
    $src = $ws.Cells.Item($r-1, 3)
    $dst = $ws.Cells.Item($r, 3)
    $dst.Value = $src.Value
    $dst.Style = $src.Style
}
$c6 = $ws.Cells.Item(6,3)
$c6.Value = "Chatter"
